# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table to the latest scraped values (GitHub Actions refresh).
#
# Numeric-looking Price values (e.g. "1.00", "137.80") are written with a
# leading apostrophe so Excel keeps them as literal text (preserving
# trailing zeros/precision) instead of auto-converting them to numbers,
# exactly like the source data (inline string cells, General format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.871.73'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').Value = '2.373.23'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').Value = '''559.51'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').Value = '''137.80'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').Value = '2.370.13'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('D12').Value = '''5.09'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').Value = '''25.67'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('E16').Value = '  -2.71%  '
$ws.Range('D17').Value = '59.809.39'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '2.371.12'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').Value = '''8.08'
$ws.Range('E19').Value = '  +13.10%  '
$ws.Range('D20').Value = '''10.51'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '''321.61'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '''4.04'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').Value = '''6.02'
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  -3.41%  '
$ws.Range('D26').Value = '''64.04'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').Value = '''556.86'
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('E28').Value = '  -6.15%  '
$ws.Range('D30').Value = '0.0₃0923'
$ws.Range('E30').Value = '  +2.02%  '
$ws.Range('D31').Value = '''8.07'
$ws.Range('E31').Value = '  +4.10%  '
$ws.Range('D32').Value = '''1.30'
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('E33').Value = '  -2.56%  '
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('E36').Value = '  +4.05%  '
$ws.Range('D37').Value = '''153.53'
$ws.Range('E37').Value = '  +4.34%  '
$ws.Range('D38').Value = '''0.366'
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').Value = '''18.14'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('D43').Value = '''41.53'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').Value = '''2.41'
$ws.Range('E45').Value = '  +2.95%  '
$ws.Range('D46').Value = '0.0₆0298'
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('D47').Value = '''139.88'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('E48').Value = '  +0.69%  '
$ws.Range('D49').Value = '''0.584'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = '''0.0499'
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('D51').Value = '''19.04'
$ws.Range('E51').Value = '  -1.30%  '
